$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: replace placeholder '?' text cells with real values / formulas ---

# B6: was shared-string "73, 3" -> now numeric 73.3
$ws.Range("B6").Value = 73.3

# C6: was shared-string "?" -> now formula =ABS(B6-73.6)
$ws.Range("C6").Formula = "=ABS(B6-73.6)"

# D6: was shared-string "?" -> now formula =ABS(B6-72.9999999999995)
$ws.Range("D6").Formula = "=ABS(B6-72.9999999999995)"

# F6: was shared-string "?" -> now formula =ABS(E6-0.314)
$ws.Range("F6").Formula = "=ABS(E6-0.314)"

# G6: was shared-string "?" -> now formula =ABS(E6-0.264)
$ws.Range("G6").Formula = "=ABS(E6-0.264)"

# Q6 / R6: new numeric inputs feeding the S6 = Q6 + R6*2 shared formula
$ws.Range("Q6").Value = 1522.7850747291
$ws.Range("R6").Value = 2

# --- Sheet view: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 92
$ws.Range("A6").Select()

# --- Column C width (Excel auto-fit "bestFit" style width) ---
# Target stored width is 11.42578125 characters; the COM layer here quantizes
# ColumnWidth to 1/6-character steps, so feed it the input that rounds to the
# closest achievable stored width (11.5).
$ws.Columns.Item(3).ColumnWidth = 10.666666666666668
